$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: fix row 10 (Objetivos:) content B10/C10 -> new Portuguese objectives text ---
$ws.Range("B10").Value = 'Discutir as consequências da natureza ondulatória da luz do ponto de vista das equações de Maxwell, familiarizando o estudante com os conceitos de interferência, difração e polarização. Apresentar os conceitos centrais da física moderna como a relatividade restrita e os fundamentos da mecânica quântica'
$ws.Range("C10").Value = 'Discutir as consequências da natureza ondulatória da luz do ponto de vista das equações de Maxwell, familiarizando o estudante com os conceitos de interferência, difração e polarização. Apresentar os conceitos centrais da física moderna como a relatividade restrita e os fundamentos da mecânica quântica'

# --- Step 2: insert a new row at 13 (shifts old rows 13-24 down to 14-25) ---
$ws.Rows.Item(13).Insert()

# Row 13 (new): carries only B/C = "230696 - Carlos Jose Todero Peixoto" (Docentes responsaveis content)
# Copy formats of B14:C14 (post-shift "Programa resumido" row, which retains the B2/C3 col style) into B13:C13,
# then clear A13 leftover + set B13/C13 text.
$ws.Range("B11:C11").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)
$ws.Range("A13").Clear()
$ws.Range("B13").Value = '230696 - Carlos José Todero Peixoto'
$ws.Range("C13").Value = '230696 - Carlos José Todero Peixoto'

# --- Step 3: row 14 (Programa resumido:) B/C -> new Portuguese short syllabus text ---
$ws.Range("B14").Value = 'Óptica Geométrica. Introdução a Física Moderna: natureza ondulatória da matéria, relatividade e introdução a mecânica quântica.'
$ws.Range("C14").Value = 'Óptica Geométrica. Introdução a Física Moderna: natureza ondulatória da matéria, relatividade e introdução a mecânica quântica.'

# --- Step 4: row 16 (Programa:) B/C -> new Portuguese full syllabus text ---
$ws.Range("B16").Value = '1) Óptica Geométrica: conceitos básicos. 2) Interferência: a experiência de Young; coerência; figuras de interferência; o interferômetro de Michelson.3) Difração.4) Polarização.5) Relatividade: os postulados da relatividade, as transformações de Lorentz, simultaneidade, tempo e comprimento; momento linear, trabalho e energia;6) Primórdios da teoria quântica: a hipótese de Plank; o efeito fotoelétrico, quantização do fóton; ondas de De Broglie, o efeito Compton, a difração de elétrons, interferência; 7) Princípios básicos da mecânica quântica: o princípio de incerteza; a equação de Schrödinger.'
$ws.Range("C16").Value = '1) Óptica Geométrica: conceitos básicos. 2) Interferência: a experiência de Young; coerência; figuras de interferência; o interferômetro de Michelson.3) Difração.4) Polarização.5) Relatividade: os postulados da relatividade, as transformações de Lorentz, simultaneidade, tempo e comprimento; momento linear, trabalho e energia;6) Primórdios da teoria quântica: a hipótese de Plank; o efeito fotoelétrico, quantização do fóton; ondas de De Broglie, o efeito Compton, a difração de elétrons, interferência; 7) Princípios básicos da mecânica quântica: o princípio de incerteza; a equação de Schrödinger.'

# --- Step 5: row 19 (Metodo:) B/C -> method text (was duplicated docentes name) ---
$ws.Range("B19").Value = 'NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.'
$ws.Range("C19").Value = 'NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.'

# --- Step 6: row 20 (Criterio:) B/C -> NF>=5,0 ---
$ws.Range("B20").Value = 'NF≥ 5,0.'
$ws.Range("C20").Value = 'NF≥ 5,0.'

# --- Step 7: row 21 (Norma de recuperacao:) B/C -> recovery formula text ---
$ws.Range("B21").Value = '(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada.'
$ws.Range("C21").Value = '(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada.'

# --- Step 8: row 22 (Bibliografia:) B/C -> bibliography text ---
$ws.Range("B22").Value = 'NUSSENZVEIG, H.M. Curso de Física Básica. Vol. 4, Edgard Blucher (2008).RESNICK, R.; HALLIDAY, D. Fundamentos de Física. Vol. 4, LTC (2008).TIPLER, P.; MOSCA, G. Física para Cientistas e Engenheiros. Vol. 4, LTC (2008).SEARS, F. W.; ZEMANSKY, M. W.; YOUNG, H. D.; FREEDMAN, R. A. Física IV, Vol. 4, Pearson Addison Wesley (2009).JEWETT Jr, John W.; SERWAY, Raymond A. Princípios de Física. Vol. 4, Thomson Pioneira (2008).'
$ws.Range("C22").Value = 'NUSSENZVEIG, H.M. Curso de Física Básica. Vol. 4, Edgard Blucher (2008).RESNICK, R.; HALLIDAY, D. Fundamentos de Física. Vol. 4, LTC (2008).TIPLER, P.; MOSCA, G. Física para Cientistas e Engenheiros. Vol. 4, LTC (2008).SEARS, F. W.; ZEMANSKY, M. W.; YOUNG, H. D.; FREEDMAN, R. A. Física IV, Vol. 4, Pearson Addison Wesley (2009).JEWETT Jr, John W.; SERWAY, Raymond A. Princípios de Física. Vol. 4, Thomson Pioneira (2008).'
